$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Shared string "2016-08-24 07:06:00" -> "2016-08-24 07:06:59"
# (used by both Overview!G2 and de-de!H2)
$wsOverview.Range("G2").Value = "2016-08-24 07:06:59"
$wsDeDe.Range("H2").Value = "2016-08-24 07:06:59"

# zh-cn!H2: "2016-08-24 07:05:55" -> "2016-08-24 07:06:53"
$wsZhCn.Range("H2").Value = "2016-08-24 07:06:53"

# zh-cn!K2: "2016-08-24 07:06:28" -> "2016-08-24 07:07:18"
$wsZhCn.Range("K2").Value = "2016-08-24 07:07:18"

# de-de!K2: "2016-08-24 07:06:36" -> "2016-08-24 07:07:26"
$wsDeDe.Range("K2").Value = "2016-08-24 07:07:26"
